$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Pros1"
$ws.Cells.Item(2, 3).Value = "Tyro3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 30.168014
$ws.Cells.Item(2, 8).Value = 90.504042
$ws.Cells.Item(2, 9).Value = 0.2349227827725553
$ws.Cells.Item(2, 10).Value = 0.2349227827725553
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.2333953333333333
$ws.Cells.Item(2, 14).Value = 0.700186
$ws.Cells.Item(2, 15).Value = 0.0679343338760815
$ws.Cells.Item(2, 16).Value = 0.0679343338760815
$ws.Cells.Item(2, 17).Value = 7.041073683534666
$ws.Cells.Item(2, 18).Value = 63.369663151812
$ws.Cells.Item(2, 19).Value = 0.01595932275996894
$ws.Cells.Item(2, 20).Value = 0.01595932275996894
# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Pros1"
$ws.Cells.Item(3, 3).Value = "Tyro3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 30.168014
$ws.Cells.Item(3, 8).Value = 90.504042
$ws.Cells.Item(3, 9).Value = 0.2349227827725553
$ws.Cells.Item(3, 10).Value = 0.2349227827725553
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.855037666666667
$ws.Cells.Item(3, 14).Value = 8.565113
$ws.Cells.Item(3, 15).Value = 0.8310152534160438
$ws.Cells.Item(3, 16).Value = 0.8310152534160438
$ws.Cells.Item(3, 17).Value = 86.13081629852734
$ws.Cells.Item(3, 18).Value = 775.177346686746
$ws.Cells.Item(3, 19).Value = 0.1952244158589373
$ws.Cells.Item(3, 20).Value = 0.1952244158589373
# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Pros1"
$ws.Cells.Item(4, 3).Value = "Tyro3"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 30.168014
$ws.Cells.Item(4, 8).Value = 90.504042
$ws.Cells.Item(4, 9).Value = 0.2349227827725553
$ws.Cells.Item(4, 10).Value = 0.2349227827725553
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.347169
$ws.Cells.Item(4, 14).Value = 1.041507
$ws.Cells.Item(4, 15).Value = 0.1010504127078748
$ws.Cells.Item(4, 16).Value = 0.1010504127078748
$ws.Cells.Item(4, 17).Value = 10.473399252366
$ws.Cells.Item(4, 18).Value = 94.260593271294
$ws.Cells.Item(4, 19).Value = 0.02373904415364913
$ws.Cells.Item(4, 20).Value = 0.02373904415364913
# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Pros1"
$ws.Cells.Item(5, 3).Value = "Tyro3"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 65.93392299999999
$ws.Cells.Item(5, 8).Value = 197.801769
$ws.Cells.Item(5, 9).Value = 0.5134372010789768
$ws.Cells.Item(5, 10).Value = 0.5134372010789768
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.2333953333333333
$ws.Cells.Item(5, 14).Value = 0.700186
$ws.Cells.Item(5, 15).Value = 0.0679343338760815
$ws.Cells.Item(5, 16).Value = 0.0679343338760815
$ws.Cells.Item(5, 17).Value = 15.38866993655933
$ws.Cells.Item(5, 18).Value = 138.498029429034
$ws.Cells.Item(5, 19).Value = 0.03488001424250001
$ws.Cells.Item(5, 20).Value = 0.03488001424250001
# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Pros1"
$ws.Cells.Item(6, 3).Value = "Tyro3"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 65.93392299999999
$ws.Cells.Item(6, 8).Value = 197.801769
$ws.Cells.Item(6, 9).Value = 0.5134372010789768
$ws.Cells.Item(6, 10).Value = 0.5134372010789768
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.855037666666667
$ws.Cells.Item(6, 14).Value = 8.565113
$ws.Cells.Item(6, 15).Value = 0.8310152534160438
$ws.Cells.Item(6, 16).Value = 0.8310152534160438
$ws.Cells.Item(6, 17).Value = 188.2438336760997
$ws.Cells.Item(6, 18).Value = 1694.194503084897
$ws.Cells.Item(6, 19).Value = 0.4266741457678702
$ws.Cells.Item(6, 20).Value = 0.4266741457678702
# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Pros1"
$ws.Cells.Item(7, 3).Value = "Tyro3"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 65.93392299999999
$ws.Cells.Item(7, 8).Value = 197.801769
$ws.Cells.Item(7, 9).Value = 0.5134372010789768
$ws.Cells.Item(7, 10).Value = 0.5134372010789768
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.347169
$ws.Cells.Item(7, 14).Value = 1.041507
$ws.Cells.Item(7, 15).Value = 0.1010504127078748
$ws.Cells.Item(7, 16).Value = 0.1010504127078748
$ws.Cells.Item(7, 17).Value = 22.890214113987
$ws.Cells.Item(7, 18).Value = 206.011927025883
$ws.Cells.Item(7, 19).Value = 0.0518830410686067
$ws.Cells.Item(7, 20).Value = 0.0518830410686067
# Row 8
$ws.Cells.Item(8, 1).Value = "M1"
$ws.Cells.Item(8, 2).Value = "Pros1"
$ws.Cells.Item(8, 3).Value = "Tyro3"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 9.984842
$ws.Cells.Item(8, 8).Value = 29.954526
$ws.Cells.Item(8, 9).Value = 0.07775344005688564
$ws.Cells.Item(8, 10).Value = 0.07775344005688566
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.2333953333333333
$ws.Cells.Item(8, 14).Value = 0.700186
$ws.Cells.Item(8, 15).Value = 0.0679343338760815
$ws.Cells.Item(8, 16).Value = 0.0679343338760815
$ws.Cells.Item(8, 17).Value = 2.330415526870667
$ws.Cells.Item(8, 18).Value = 20.973739741836
$ws.Cells.Item(8, 19).Value = 0.005282128156838359
$ws.Cells.Item(8, 20).Value = 0.00528212815683836
# Row 9
$ws.Cells.Item(9, 1).Value = "M1"
$ws.Cells.Item(9, 2).Value = "Pros1"
$ws.Cells.Item(9, 3).Value = "Tyro3"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 9.984842
$ws.Cells.Item(9, 8).Value = 29.954526
$ws.Cells.Item(9, 9).Value = 0.07775344005688564
$ws.Cells.Item(9, 10).Value = 0.07775344005688566
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.855037666666667
$ws.Cells.Item(9, 14).Value = 8.565113
$ws.Cells.Item(9, 15).Value = 0.8310152534160438
$ws.Cells.Item(9, 16).Value = 0.8310152534160438
$ws.Cells.Item(9, 17).Value = 28.50710000571533
$ws.Cells.Item(9, 18).Value = 256.563900051438
$ws.Cells.Item(9, 19).Value = 0.06461429469284199
$ws.Cells.Item(9, 20).Value = 0.06461429469284199
# Row 10
$ws.Cells.Item(10, 1).Value = "M1"
$ws.Cells.Item(10, 2).Value = "Pros1"
$ws.Cells.Item(10, 3).Value = "Tyro3"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 9.984842
$ws.Cells.Item(10, 8).Value = 29.954526
$ws.Cells.Item(10, 9).Value = 0.07775344005688564
$ws.Cells.Item(10, 10).Value = 0.07775344005688566
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.347169
$ws.Cells.Item(10, 14).Value = 1.041507
$ws.Cells.Item(10, 15).Value = 0.1010504127078748
$ws.Cells.Item(10, 16).Value = 0.1010504127078748
$ws.Cells.Item(10, 17).Value = 3.466427612298
$ws.Cells.Item(10, 18).Value = 31.197848510682
$ws.Cells.Item(10, 19).Value = 0.007857017207205297
$ws.Cells.Item(10, 20).Value = 0.007857017207205299
# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Pros1"
$ws.Cells.Item(11, 3).Value = "Tyro3"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 16.22962566666667
$ws.Cells.Item(11, 8).Value = 48.688877
$ws.Cells.Item(11, 9).Value = 0.1263824932251166
$ws.Cells.Item(11, 10).Value = 0.1263824932251166
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.2333953333333333
$ws.Cells.Item(11, 14).Value = 0.700186
$ws.Cells.Item(11, 15).Value = 0.0679343338760815
$ws.Cells.Item(11, 16).Value = 0.0679343338760815
$ws.Cells.Item(11, 17).Value = 3.787918892346889
$ws.Cells.Item(11, 18).Value = 34.091270031122
$ws.Cells.Item(11, 19).Value = 0.008585710490846677
$ws.Cells.Item(11, 20).Value = 0.008585710490846677
# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Pros1"
$ws.Cells.Item(12, 3).Value = "Tyro3"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 16.22962566666667
$ws.Cells.Item(12, 8).Value = 48.688877
$ws.Cells.Item(12, 9).Value = 0.1263824932251166
$ws.Cells.Item(12, 10).Value = 0.1263824932251166
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 2.855037666666667
$ws.Cells.Item(12, 14).Value = 8.565113
$ws.Cells.Item(12, 15).Value = 0.8310152534160438
$ws.Cells.Item(12, 16).Value = 0.8310152534160438
$ws.Cells.Item(12, 17).Value = 46.33619259423345
$ws.Cells.Item(12, 18).Value = 417.025733348101
$ws.Cells.Item(12, 19).Value = 0.1050257796348217
$ws.Cells.Item(12, 20).Value = 0.1050257796348217
# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Pros1"
$ws.Cells.Item(13, 3).Value = "Tyro3"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 16.22962566666667
$ws.Cells.Item(13, 8).Value = 48.688877
$ws.Cells.Item(13, 9).Value = 0.1263824932251166
$ws.Cells.Item(13, 10).Value = 0.1263824932251166
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.347169
$ws.Cells.Item(13, 14).Value = 1.041507
$ws.Cells.Item(13, 15).Value = 0.1010504127078748
$ws.Cells.Item(13, 16).Value = 0.1010504127078748
$ws.Cells.Item(13, 17).Value = 5.634422913071
$ws.Cells.Item(13, 18).Value = 50.709806217639
$ws.Cells.Item(13, 19).Value = 0.01277100309944822
$ws.Cells.Item(13, 20).Value = 0.01277100309944822
# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Pros1"
$ws.Cells.Item(14, 3).Value = "Tyro3"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 6.100318666666666
$ws.Cells.Item(14, 8).Value = 18.300956
$ws.Cells.Item(14, 9).Value = 0.04750408286646571
$ws.Cells.Item(14, 10).Value = 0.04750408286646571
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 0.2333953333333333
$ws.Cells.Item(14, 14).Value = 0.700186
$ws.Cells.Item(14, 15).Value = 0.0679343338760815
$ws.Cells.Item(14, 16).Value = 0.0679343338760815
$ws.Cells.Item(14, 17).Value = 1.423785908646222
$ws.Cells.Item(14, 18).Value = 12.814073177816
$ws.Cells.Item(14, 19).Value = 0.003227158225927524
$ws.Cells.Item(14, 20).Value = 0.003227158225927524
# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Pros1"
$ws.Cells.Item(15, 3).Value = "Tyro3"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 6.100318666666666
$ws.Cells.Item(15, 8).Value = 18.300956
$ws.Cells.Item(15, 9).Value = 0.04750408286646571
$ws.Cells.Item(15, 10).Value = 0.04750408286646571
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 2.855037666666667
$ws.Cells.Item(15, 14).Value = 8.565113
$ws.Cells.Item(15, 15).Value = 0.8310152534160438
$ws.Cells.Item(15, 16).Value = 0.8310152534160438
$ws.Cells.Item(15, 17).Value = 17.41663957200311
$ws.Cells.Item(15, 18).Value = 156.749756148028
$ws.Cells.Item(15, 19).Value = 0.03947661746157274
$ws.Cells.Item(15, 20).Value = 0.03947661746157275
# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Pros1"
$ws.Cells.Item(16, 3).Value = "Tyro3"
$ws.Cells.Item(16, 4).Value = "sCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 6.100318666666666
$ws.Cells.Item(16, 8).Value = 18.300956
$ws.Cells.Item(16, 9).Value = 0.04750408286646571
$ws.Cells.Item(16, 10).Value = 0.04750408286646571
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.347169
$ws.Cells.Item(16, 14).Value = 1.041507
$ws.Cells.Item(16, 15).Value = 0.1010504127078748
$ws.Cells.Item(16, 16).Value = 0.1010504127078748
$ws.Cells.Item(16, 17).Value = 2.117841531188
$ws.Cells.Item(16, 18).Value = 19.060573780692
$ws.Cells.Item(16, 19).Value = 0.004800307178965442
$ws.Cells.Item(16, 20).Value = 0.004800307178965443
